$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row (column names)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Title-case state/municipality names that contain Spanish connector words
# (de/del/el/la/las/los/y) so they read "De"/"Del"/"El"/"La"/"Las"/"Los"/"Y"
$ws.Range("B21").Value = "Amatenango De La Frontera"
$ws.Range("B25").Value = "Bejucal De Ocampo"
$ws.Range("B36").Value = "Comitán De Domínguez"
$ws.Range("B59").Value = "Mazapa De Madero"
$ws.Range("B62").Value = "Montecristo De Guerrero"
$ws.Range("B65").Value = "Ocozocoautla De Espinosa"
$ws.Range("B70").Value = "Salto De Agua"
$ws.Range("B71").Value = "San Cristóbal De Las Casas"
$ws.Range("B95").Value = "Hidalgo Del Parral"
$ws.Range("B100").Value = "San Francisco Del Oro"
$ws.Range("A103").Value = "Ciudad De México"
$ws.Range("B106").Value = "Cuajimalpa De Morelos"
$ws.Range("A121").Value = "Coahuila De Zaragoza"
$ws.Range("B133").Value = "San Juan Del Río"
$ws.Range("A135").Value = "Estado De México"
$ws.Range("B135").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B137").Value = "Almoloya De Alquisiras"
$ws.Range("B138").Value = "Almoloya De Juárez"
$ws.Range("B139").Value = "Almoloya Del Río"
$ws.Range("B144").Value = "Atizapán De Zaragoza"
$ws.Range("B150").Value = "Chapa De Mota"
$ws.Range("B159").Value = "Ecatepec De Morelos"
$ws.Range("B164").Value = "Ixtapan De La Sal"
$ws.Range("B165").Value = "Ixtapan Del Oro"
$ws.Range("B173").Value = "Naucalpan De Juárez"
$ws.Range("B180").Value = "San Felipe Del Progreso"
$ws.Range("B181").Value = "San José Del Rincón"
$ws.Range("B183").Value = "San Simón De Guerrero"
$ws.Range("B191").Value = "Tenango Del Valle"
$ws.Range("B198").Value = "Tlalnepantla De Baz"
$ws.Range("B202").Value = "Valle De Bravo"
$ws.Range("B203").Value = "Valle De Chalco Solidaridad"
$ws.Range("B206").Value = "Villa De Allende"
$ws.Range("B216").Value = "Apaseo El Alto"
$ws.Range("B217").Value = "Apaseo El Grande"
$ws.Range("B223").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B227").Value = "Jaral Del Progreso"
$ws.Range("B235").Value = "San Diego De La Unión"
$ws.Range("B237").Value = "San Francisco Del Rincón"
$ws.Range("B238").Value = "San Miguel De Allende"
$ws.Range("B239").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B241").Value = "Silao De La Victoria"
$ws.Range("B244").Value = "Valle De Santiago"
$ws.Range("B250").Value = "Acapulco De Juárez"
$ws.Range("B253").Value = "Ajuchitlán Del Progreso"
$ws.Range("B254").Value = "Alcozauca De Guerrero"
$ws.Range("B257").Value = "Atlamajalcingo Del Monte"
$ws.Range("B259").Value = "Atoyac De Álvarez"
$ws.Range("B260").Value = "Ayutla De Los Libres"
$ws.Range("B262").Value = "Chilapa De Álvarez"
$ws.Range("B263").Value = "Chilpancingo De Los Bravo"
$ws.Range("B264").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B268").Value = "Coyuca De Benítez"
$ws.Range("B269").Value = "Coyuca De Catalán"
$ws.Range("B272").Value = "Cuetzala Del Progreso"
$ws.Range("B273").Value = "Cutzamala De Pinzón"
$ws.Range("B280").Value = "Huitzuco De Los Figueroa"
$ws.Range("B281").Value = "Iguala De La Independencia"
$ws.Range("B282").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B297").Value = "Taxco De Alarcón"
$ws.Range("B300").Value = "Tepecoacuilco De Trujano"
$ws.Range("B301").Value = "Tixtla De Guerrero"
$ws.Range("B304").Value = "Tlapa De Comonfort"
$ws.Range("B305").Value = "Técpan De Galeana"
$ws.Range("B310").Value = "Zihuatanejo De Azueta"
$ws.Range("B318").Value = "Atotonilco El Grande"
$ws.Range("B322").Value = "Cuautepec De Hinojosa"
$ws.Range("B324").Value = "Huasca De Ocampo"
$ws.Range("B331").Value = "Mineral Del Chico"
$ws.Range("B332").Value = "Mixquiahuala De Juárez"
$ws.Range("B334").Value = "Omitlán De Juárez"
$ws.Range("B335").Value = "Pachuca De Soto"
$ws.Range("B337").Value = "Progreso De Obregón"
$ws.Range("B343").Value = "Santiago De Anaya"
$ws.Range("B346").Value = "Tenango De Doria"
$ws.Range("B348").Value = "Tepehuacán De Guerrero"
$ws.Range("B349").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B350").Value = "Tezontepec De Aldama"
$ws.Range("B354").Value = "Tulancingo De Bravo"
$ws.Range("B357").Value = "Zacualtipán De Ángeles"
$ws.Range("B358").Value = "Zapotlán De Juárez"
$ws.Range("B361").Value = "Ahualulco De Mercado"
$ws.Range("B362").Value = "Atotonilco El Alto"
$ws.Range("B366").Value = "Encarnación De Díaz"
$ws.Range("B369").Value = "Huejuquilla El Alto"
$ws.Range("B370").Value = "Ixtlahuacán De Los Membrillos"
$ws.Range("B375").Value = "Lagos De Moreno"
$ws.Range("B378").Value = "Ojuelos De Jalisco"
$ws.Range("B379").Value = "San Juan De Los Lagos"
$ws.Range("B381").Value = "San Martín De Bolaños"
$ws.Range("A389").Value = "Michoacán De Ocampo"
$ws.Range("B441").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B457").Value = "Coatlán Del Río"
$ws.Range("B468").Value = "Puente De Ixtla"
$ws.Range("B470").Value = "Tetela Del Volcán"
$ws.Range("B490").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B493").Value = "Chiquihuitlán De Benito Juárez"
$ws.Range("B496").Value = "Cuilápam De Guerrero"
$ws.Range("B498").Value = "Guevea De Humboldt"
$ws.Range("B499").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B500").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B501").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B502").Value = "Huajuapan De León"
$ws.Range("B506").Value = "Mazatlán Villa De Flores"
$ws.Range("B508").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B510").Value = "Oaxaca De Juárez"
$ws.Range("B511").Value = "Ocotlán De Morelos"
$ws.Range("B512").Value = "Putla Villa De Guerrero"
$ws.Range("B524").Value = "San Dionisio Del Mar"
$ws.Range("B528").Value = "San Francisco Del Mar"
$ws.Range("B531").Value = "San José Del Progreso"
$ws.Range("B533").Value = "San Juan Bautista Lo De Soto"
$ws.Range("B553").Value = "San Pablo Villa De Mitla"
$ws.Range("B558").Value = "San Pedro Y San Pablo Ayutla"
$ws.Range("B562").Value = "Santa Inés De Zaragoza"
$ws.Range("B590").Value = "Santo Domingo De Morelos"
$ws.Range("B594").Value = "Tataltepec De Valdés"
$ws.Range("B595").Value = "Teotitlán De Flores Magón"
$ws.Range("B596").Value = "Villa De Chilapa De Díaz"
$ws.Range("B597").Value = "Villa De Etla"
$ws.Range("B598").Value = "Villa De Tututepec"
$ws.Range("B599").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B601").Value = "Zimatlán De Álvarez"
$ws.Range("B629").Value = "Cuapiaxtla De Madero"
$ws.Range("B631").Value = "Cuayuca De Andrade"
$ws.Range("B632").Value = "Cuetzalan Del Progreso"
$ws.Range("B645").Value = "Huehuetlán El Chico"
$ws.Range("B646").Value = "Huehuetlán El Grande"
$ws.Range("B649").Value = "Huitzilan De Serdán"
$ws.Range("B651").Value = "Ixcamilpa De Guerrero"
$ws.Range("B654").Value = "Izúcar De Matamoros"
$ws.Range("B662").Value = "Los Reyes De Juárez"
$ws.Range("B671").Value = "Palmar De Bravo"
$ws.Range("B688").Value = "San Nicolás De Los Ranchos"
$ws.Range("B692").Value = "San Salvador El Seco"
$ws.Range("B693").Value = "San Salvador El Verde"
$ws.Range("B696").Value = "Tecali De Herrera"
$ws.Range("B703").Value = "Tepanco De López"
$ws.Range("B709").Value = "Tepexi De Rodríguez"
$ws.Range("B711").Value = "Tepeyahualco De Cuauhtémoc"
$ws.Range("B712").Value = "Tetela De Ocampo"
$ws.Range("B716").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B725").Value = "Tuzamapan De Galeana"
$ws.Range("B731").Value = "Xochitlán De Vicente Suárez"
$ws.Range("B741").Value = "Amealco De Bonfil"
$ws.Range("B742").Value = "Cadereyta De Montes"
$ws.Range("B747").Value = "Jalpan De Serra"
$ws.Range("B749").Value = "Pinal De Amoles"
$ws.Range("B752").Value = "San Juan Del Río"
$ws.Range("B763").Value = "Axtla De Terrazas"
$ws.Range("B771").Value = "Mexquitic De Carmona"
$ws.Range("B777").Value = "Santa María Del Río"
$ws.Range("B781").Value = "Villa De Guadalupe"
$ws.Range("B802").Value = "Jalpa De Méndez"
$ws.Range("B812").Value = "Soto La Marina"
$ws.Range("B816").Value = "Apetatitlán De Antonio Carvajal"
$ws.Range("B828").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B830").Value = "Nanacamilpa De Mariano Arista"
$ws.Range("B836").Value = "Tepetitla De Lardizábal"
$ws.Range("A846").Value = "Veracruz De Ignacio De La Llave"
$ws.Range("B854").Value = "Boca Del Río"
$ws.Range("B863").Value = "Cosamaloapan De Carpio"
$ws.Range("B875").Value = "Hueyapan De Ocampo"
$ws.Range("B876").Value = "Ignacio De La Llave"
$ws.Range("B879").Value = "Ixhuatlán De Madero"
$ws.Range("B885").Value = "Lerdo De Tejada"
$ws.Range("B887").Value = "Martínez De La Torre"
$ws.Range("B893").Value = "Mixtla De Altamirano"
$ws.Range("B903").Value = "Poza Rica De Hidalgo"
$ws.Range("B908").Value = "Sayula De Alemán"
$ws.Range("B910").Value = "Soledad De Doblado"
$ws.Range("B932").Value = "Vega De Alatorre"
$ws.Range("B939").Value = "Zontecomatlán De López Y Fuentes"
$ws.Range("A956").Value = "Total"

# Remove the trailing footnote rows (sample size / source / author / date)
$ws.Rows("959:963").Delete()
